$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- "coloring" requirement block (rows 5-7) ---
# Row5: world model access wording -> world model action wording
$ws.Range("D5").Value = "world model action을 coloring 할 수 있어야 한다."
$ws.Range("E5").Value = "fact, retrieve, match, update, retract, assert 단어를 coloring 할 수 있어야 한다."

# Row6 used to hold the "prefix" requirement; it now holds the
# "특정 keyword / plan, import" requirement (moved up from old row7),
# fully scored (중요도5, 복잡도1, 완료 o).
$ws.Range("D6").Value = "특정 keyword를 coloring 할 수 있어야 한다."
$ws.Range("E6").Value = "plan, import"
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = "o"

# Row7 now holds the "prefix" requirement (moved down from old row6).
$ws.Range("D7").Value = "prefix 단어를 coloring 할 수 있어야 한다."
$ws.Range("E7").Value = "prefixes 안의 prefix된 단어들을 coloring 할 수 있어야 한다."
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1

# Row8 ("assert, retract도 해야 함") is removed entirely - that remark is now folded into row5's wording.
$ws.Range("D8").ClearContents()

# --- 자동완성 (auto-complete) block (rows 10-13) ---
$ws.Range("F10").Value = 1
$ws.Range("H10").Value = "x"
$ws.Range("I10").Value = "없애도 될 듯"

$ws.Range("H11").Value = "o"
$ws.Range("I11").Value = "body:가 아닌 빈 공간에서 자동 완성 커맨드(Ctrl + space)사용 시 "

$ws.Range("H12").Value = "o"
$ws.Range("I12").Value = "body:가 goal action 이후 자동 완성 커맨드(Ctrl + space)사용 시"

$ws.Range("H13").Value = "o"

# --- 문법 오류 찾기 (syntax-error) block (rows 19-23) now all marked complete ---
$ws.Range("H19").Value = "o"
$ws.Range("H20").Value = "o"
$ws.Range("H21").Value = "o"
$ws.Range("H22").Value = "o"
$ws.Range("H23").Value = "o"

# --- resource hyperlink block (row 27) ---
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = "o"

# --- jam plan model block (rows 29-30) ---
$ws.Range("H29").Value = "o"
$ws.Range("H30").Value = "o"

# Match style of existing "완료" (H) column cells (centered, like H2:H5).
$ws.Range("H6,H10,H11,H12,H13,H19,H20,H21,H22,H23,H27,H29,H30").HorizontalAlignment = -4108

$ws.Range("F8").Select()

$wb.Save()
